$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Barney"
$ws.Range("B4").Value = "Fife"
$ws.Range("C4").Value = "333-11-2345"

$ws.Range("C5").Select()
